$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text such as "23.481.06" or "1.004" that must
# stay literal text (Excel would otherwise coerce them into numbers and drop
# formatting like trailing zeros or the multi-dot "thousands" grouping used
# on this sheet), so force a text number format before writing the values.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.481.06'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.646.13'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '302.63'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3841'
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3595'
$ws.Range('E8').Value = '  +1.15%  '
$ws.Range('E9').Value = '  -1.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08169'
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.229'
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.004'
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.30'
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.444'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.458'
$ws.Range('E15').Value = '  +2.12%  '
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.654.34'
$ws.Range('E17').Value = '  +1.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '97.46'
$ws.Range('E18').Value = '  +2.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07008'
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('E20').Value = '  +2.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.54'
$ws.Range('E21').Value = '  +1.28%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.60'
$ws.Range('E23').Value = '  +1.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.491.41'
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('E25').Value = '  -2.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.033'
$ws.Range('E26').Value = '  -2.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.20'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '152.93'
$ws.Range('E28').Value = '  +1.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.226'
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.91'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.832.62'
$ws.Range('E31').Value = '  +1.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.115'
$ws.Range('E32').Value = '  +9.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.246'
$ws.Range('E33').Value = '  +4.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.19'
$ws.Range('E34').Value = '  +5.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.054'
$ws.Range('E35').Value = '  -1.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02789'
$ws.Range('E36').Value = '  +1.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2501'
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.067'
$ws.Range('E39').Value = '  +2.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06978'
$ws.Range('E40').Value = '  +0.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.05'
$ws.Range('E41').Value = '  +7.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6971'
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.333'
$ws.Range('E43').Value = '  +0.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.88'
$ws.Range('E44').Value = '  +3.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6500'
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.003'
$ws.Range('E46').Value = '  +0.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.293'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.954'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07869'
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '128.06'
$ws.Range('E50').Value = '  -1.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.176'
$ws.Range('E51').Value = '  -0.54%  '
